# Generate Report for Handback
# Update the "generated at" timestamps recorded in the handback status report.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
# "Latest HO Xliff Generate Date" for the first row moves forward.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-24 09:08:22"

# --- zh-cn sheet ---
# "Correspond Handoff Datetime" / "Correspond Handback DateTime" for the first row move forward.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-24 09:08:17"
$wsZhCn.Range("K2").Value = "2016-08-24 09:08:35"

# --- de-de sheet ---
# "Correspond Handoff Datetime" (H2) holds the same timestamp text as Overview!G2,
# and "Correspond Handback DateTime" (K2) for the first row moves forward.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-24 09:08:22"
$wsDeDe.Range("K2").Value = "2016-08-24 09:08:43"
